$wb = $excel.ActiveWorkbook

# --- ALC (sheet index 1) ---
$ws = $wb.Worksheets.Item(1)
# Row 98
$ws.Cells.Item(98, 8).Value = 2680.3333
$ws.Cells.Item(98, 9).Value = 2622.9565
$ws.Cells.Item(98, 11).Value = 2622.9565
$ws.Cells.Item(98, 13).Value = -1124.9565
# Row 122
$ws.Cells.Item(122, 8).Value = 2680.3333
$ws.Cells.Item(122, 9).Value = 2622.9565
$ws.Cells.Item(122, 11).Value = 7868.869499999999
$ws.Cells.Item(122, 13).Value = -5418.869499999999
# Row 137
$ws.Cells.Item(137, 8).Value = 18183008
$ws.Cells.Item(137, 9).Value = 1192.9688
$ws.Cells.Item(137, 11).Value = 3578.9064
$ws.Cells.Item(137, 13).Value = -1028.9064
# Row 138
$ws.Cells.Item(138, 8).Value = 2218.75
$ws.Cells.Item(138, 9).Value = 1644.2046
$ws.Cells.Item(138, 10).Value = 3121.6072
$ws.Cells.Item(138, 11).Value = 4932.6138
$ws.Cells.Item(138, 12).Value = 9364.821599999999
$ws.Cells.Item(138, 13).Value = 207.3861999999999
$ws.Cells.Item(138, 14).Value = -19644.8216
# Row 141
$ws.Cells.Item(141, 8).Value = 1106.5577
$ws.Cells.Item(141, 9).Value = 680.36365
$ws.Cells.Item(141, 11).Value = 2041.09095
$ws.Cells.Item(141, 13).Value = 3138.90905

# --- ARM (sheet index 2) ---
$ws = $wb.Worksheets.Item(2)
# Row 2
$ws.Cells.Item(2, 8).Value = 920245.4
$ws.Cells.Item(2, 9).Value = 902.2222
$ws.Cells.Item(2, 10).Value = 2102258
$ws.Cells.Item(2, 11).Value = 902.2222
$ws.Cells.Item(2, 12).Value = 2102258
$ws.Cells.Item(2, 13).Value = -789.2222
$ws.Cells.Item(2, 14).Value = -2102484
# Row 32
$ws.Cells.Item(32, 8).Value = 6172.093
$ws.Cells.Item(32, 9).Value = 4485.9287
$ws.Cells.Item(32, 10).Value = 17067.309
$ws.Cells.Item(32, 11).Value = 4485.9287
$ws.Cells.Item(32, 12).Value = 17067.309
$ws.Cells.Item(32, 13).Value = -4198.9287
$ws.Cells.Item(32, 14).Value = -17641.309
# Row 110
$ws.Cells.Item(110, 8).Value = 1554.4
$ws.Cells.Item(110, 9).Value = 1220.5
$ws.Cells.Item(110, 10).Value = 2890
$ws.Cells.Item(110, 11).Value = 1220.5
$ws.Cells.Item(110, 12).Value = 2890
$ws.Cells.Item(110, 13).Value = 824.5
$ws.Cells.Item(110, 14).Value = -6980
# Row 116
$ws.Cells.Item(116, 8).Value = 920245.4
$ws.Cells.Item(116, 9).Value = 902.2222
$ws.Cells.Item(116, 10).Value = 2102258
$ws.Cells.Item(116, 11).Value = 902.2222
$ws.Cells.Item(116, 12).Value = 2102258
$ws.Cells.Item(116, 13).Value = 1391.7778
$ws.Cells.Item(116, 14).Value = -2106846
# Row 122
$ws.Cells.Item(122, 8).Value = 1346.0303
$ws.Cells.Item(122, 9).Value = 1300.826
$ws.Cells.Item(122, 10).Value = 1450
$ws.Cells.Item(122, 11).Value = 3902.478
$ws.Cells.Item(122, 12).Value = 4350
$ws.Cells.Item(122, 13).Value = -1452.478
$ws.Cells.Item(122, 14).Value = -9250

# --- BSM (sheet index 3) ---
$ws = $wb.Worksheets.Item(3)
# Row 3
$ws.Cells.Item(3, 8).Value = 920245.4
$ws.Cells.Item(3, 9).Value = 902.2222
$ws.Cells.Item(3, 10).Value = 2102258
$ws.Cells.Item(3, 11).Value = 902.2222
$ws.Cells.Item(3, 12).Value = 2102258
$ws.Cells.Item(3, 13).Value = -788.2222
$ws.Cells.Item(3, 14).Value = -2102486
# Row 94
$ws.Cells.Item(94, 8).Value = 748.65216
$ws.Cells.Item(94, 9).Value = 658.8946999999999
$ws.Cells.Item(94, 10).Value = 1175
$ws.Cells.Item(94, 11).Value = 658.8946999999999
$ws.Cells.Item(94, 12).Value = 1175
$ws.Cells.Item(94, 13).Value = -207.8946999999999
$ws.Cells.Item(94, 14).Value = -2077
# Row 134
$ws.Cells.Item(134, 8).Value = 4469723.5
$ws.Cells.Item(134, 9).Value = 4905662
$ws.Cells.Item(134, 10).Value = 1350.25
$ws.Cells.Item(134, 11).Value = 14716986
$ws.Cells.Item(134, 12).Value = 4050.75
$ws.Cells.Item(134, 13).Value = -14714451
$ws.Cells.Item(134, 14).Value = -9120.75

# --- CRP (sheet index 4) ---
$ws = $wb.Worksheets.Item(4)
# Row 22
$ws.Cells.Item(22, 8).Value = 201.21875
$ws.Cells.Item(22, 9).Value = 193.11111
$ws.Cells.Item(22, 10).Value = 245
$ws.Cells.Item(22, 11).Value = 193.11111
$ws.Cells.Item(22, 12).Value = 245
$ws.Cells.Item(22, 13).Value = 156.88889
$ws.Cells.Item(22, 14).Value = -945
# Row 31
$ws.Cells.Item(31, 8).Value = 5594718
$ws.Cells.Item(31, 9).Value = 1517.2
$ws.Cells.Item(31, 10).Value = 10255718
$ws.Cells.Item(31, 11).Value = 1517.2
$ws.Cells.Item(31, 12).Value = 10255718
$ws.Cells.Item(31, 13).Value = -1222.2
$ws.Cells.Item(31, 14).Value = -10256308
# Row 34
$ws.Cells.Item(34, 8).Value = 5594718
$ws.Cells.Item(34, 9).Value = 1517.2
$ws.Cells.Item(34, 10).Value = 10255718
$ws.Cells.Item(34, 11).Value = 1517.2
$ws.Cells.Item(34, 12).Value = 10255718
$ws.Cells.Item(34, 13).Value = -1315.2
$ws.Cells.Item(34, 14).Value = -10256122
# Row 53
$ws.Cells.Item(53, 8).Value = 35000
$ws.Cells.Item(53, 10).Value = 35000
$ws.Cells.Item(53, 12).Value = 35000
$ws.Cells.Item(53, 14).Value = -36214
# Row 99
$ws.Cells.Item(99, 8).Value = 2279.1765
$ws.Cells.Item(99, 9).Value = 2115.4666
$ws.Cells.Item(99, 10).Value = 3507
$ws.Cells.Item(99, 11).Value = 2115.4666
$ws.Cells.Item(99, 12).Value = 3507
$ws.Cells.Item(99, 13).Value = -617.4666000000002
$ws.Cells.Item(99, 14).Value = -6503
# Row 111
$ws.Cells.Item(111, 8).Value = 59620.8
$ws.Cells.Item(111, 10).Value = 59620.8
$ws.Cells.Item(111, 12).Value = 59620.8
$ws.Cells.Item(111, 14).Value = -67800.8
# Row 122
$ws.Cells.Item(122, 8).Value = 1305.3
$ws.Cells.Item(122, 9).Value = 1150.5714
$ws.Cells.Item(122, 10).Value = 1666.3334
$ws.Cells.Item(122, 11).Value = 3451.7142
$ws.Cells.Item(122, 12).Value = 4999.0002
$ws.Cells.Item(122, 13).Value = -1001.7142
$ws.Cells.Item(122, 14).Value = -9899.0002
# Row 126
$ws.Cells.Item(126, 8).Value = 2279.1765
$ws.Cells.Item(126, 9).Value = 2115.4666
$ws.Cells.Item(126, 10).Value = 3507
$ws.Cells.Item(126, 11).Value = 6346.399800000001
$ws.Cells.Item(126, 12).Value = 10521
$ws.Cells.Item(126, 13).Value = -3876.399800000001
$ws.Cells.Item(126, 14).Value = -15461

# --- CUL (sheet index 5) ---
$ws = $wb.Worksheets.Item(5)
# Row 51
$ws.Cells.Item(51, 8).Value = 1939.2858
$ws.Cells.Item(51, 10).Value = 1939.2858
$ws.Cells.Item(51, 12).Value = 5817.857400000001
$ws.Cells.Item(51, 14).Value = -6737.857400000001
# Row 68
$ws.Cells.Item(68, 8).Value = 1571.1515
$ws.Cells.Item(68, 9).Value = 834.6667
$ws.Cells.Item(68, 10).Value = 1847.3334
$ws.Cells.Item(68, 11).Value = 2504.0001
$ws.Cells.Item(68, 12).Value = 5542.0002
$ws.Cells.Item(68, 13).Value = -1693.0001
$ws.Cells.Item(68, 14).Value = -7164.0002
# Row 71
$ws.Cells.Item(71, 8).Value = 1571.1515
$ws.Cells.Item(71, 9).Value = 834.6667
$ws.Cells.Item(71, 10).Value = 1847.3334
$ws.Cells.Item(71, 11).Value = 7512.0003
$ws.Cells.Item(71, 12).Value = 16626.0006
$ws.Cells.Item(71, 13).Value = -3456.0003
$ws.Cells.Item(71, 14).Value = -24738.0006
# Row 110
$ws.Cells.Item(110, 8).Value = 2643.3333
$ws.Cells.Item(110, 10).Value = 3030
$ws.Cells.Item(110, 12).Value = 9090
$ws.Cells.Item(110, 14).Value = -17270
# Row 131
$ws.Cells.Item(131, 8).Value = 3550.05
$ws.Cells.Item(131, 9).Value = 5469.1665
$ws.Cells.Item(131, 10).Value = 2727.5715
$ws.Cells.Item(131, 11).Value = 16407.4995
$ws.Cells.Item(131, 12).Value = 8182.7145
$ws.Cells.Item(131, 13).Value = -11367.4995
$ws.Cells.Item(131, 14).Value = -18262.7145

# --- GSM (sheet index 6) ---
$ws = $wb.Worksheets.Item(6)
# Row 46
$ws.Cells.Item(46, 8).Value = 0
$ws.Cells.Item(46, 10).Value = 0
$ws.Cells.Item(46, 12).Value = 0
$ws.Cells.Item(46, 14).ClearContents()

# --- LTW (sheet index 7) ---
$ws = $wb.Worksheets.Item(7)
# Row 7
$ws.Cells.Item(7, 8).Value = 4490.4
$ws.Cells.Item(7, 9).Value = 4150.6665
$ws.Cells.Item(7, 10).Value = 5000
$ws.Cells.Item(7, 11).Value = 4150.6665
$ws.Cells.Item(7, 12).Value = 5000
$ws.Cells.Item(7, 13).Value = -4038.6665
$ws.Cells.Item(7, 14).Value = -5224
# Row 40
$ws.Cells.Item(40, 8).Value = 1599.7
$ws.Cells.Item(40, 9).Value = 1599.7
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 11).Value = 1599.7
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 13).Value = -1463.7
$ws.Cells.Item(40, 14).ClearContents()
# Row 82
$ws.Cells.Item(82, 8).Value = 1411.579
$ws.Cells.Item(82, 9).Value = 1221.8182
$ws.Cells.Item(82, 10).Value = 1672.5
$ws.Cells.Item(82, 11).Value = 1221.8182
$ws.Cells.Item(82, 12).Value = 1672.5
$ws.Cells.Item(82, 13).Value = -860.8181999999999
$ws.Cells.Item(82, 14).Value = -2394.5
# Row 85
$ws.Cells.Item(85, 8).Value = 1411.579
$ws.Cells.Item(85, 9).Value = 1221.8182
$ws.Cells.Item(85, 10).Value = 1672.5
$ws.Cells.Item(85, 11).Value = 1221.8182
$ws.Cells.Item(85, 12).Value = 1672.5
$ws.Cells.Item(85, 13).Value = 26.18180000000007
$ws.Cells.Item(85, 14).Value = -4168.5
# Row 100
$ws.Cells.Item(100, 8).Value = 3020.2
$ws.Cells.Item(100, 9).Value = 1501
$ws.Cells.Item(100, 10).Value = 4033
$ws.Cells.Item(100, 11).Value = 1501
$ws.Cells.Item(100, 12).Value = 4033
$ws.Cells.Item(100, 13).Value = -960
$ws.Cells.Item(100, 14).Value = -5115
# Row 122
$ws.Cells.Item(122, 8).Value = 1822.8572
$ws.Cells.Item(122, 9).Value = 1738.1818
$ws.Cells.Item(122, 11).Value = 5214.5454
$ws.Cells.Item(122, 13).Value = -2764.5454
# Row 126
$ws.Cells.Item(126, 8).Value = 4490.4
$ws.Cells.Item(126, 9).Value = 4150.6665
$ws.Cells.Item(126, 10).Value = 5000
$ws.Cells.Item(126, 11).Value = 12451.9995
$ws.Cells.Item(126, 12).Value = 15000
$ws.Cells.Item(126, 13).Value = -9981.999500000002
$ws.Cells.Item(126, 14).Value = -19940
# Row 132
$ws.Cells.Item(132, 8).Value = 3616
$ws.Cells.Item(132, 9).Value = 4037.342
$ws.Cells.Item(132, 10).Value = 413.8
$ws.Cells.Item(132, 11).Value = 12112.026
$ws.Cells.Item(132, 12).Value = 1241.4
$ws.Cells.Item(132, 13).Value = -9582.026
$ws.Cells.Item(132, 14).Value = -6301.4
# Row 136
$ws.Cells.Item(136, 8).Value = 1488.5834
$ws.Cells.Item(136, 9).Value = 737
$ws.Cells.Item(136, 10).Value = 3196.7273
$ws.Cells.Item(136, 11).Value = 2211
$ws.Cells.Item(136, 12).Value = 9590.1819
$ws.Cells.Item(136, 13).Value = 339
$ws.Cells.Item(136, 14).Value = -14690.1819

# --- WVR (sheet index 8) ---
$ws = $wb.Worksheets.Item(8)
# Row 80
$ws.Cells.Item(80, 8).Value = 10000
$ws.Cells.Item(80, 10).Value = 10000
$ws.Cells.Item(80, 12).Value = 10000
$ws.Cells.Item(80, 14).Value = -11996
# Row 83
$ws.Cells.Item(83, 8).Value = 10000
$ws.Cells.Item(83, 10).Value = 10000
$ws.Cells.Item(83, 12).Value = 30000
$ws.Cells.Item(83, 14).Value = -39984
# Row 132
$ws.Cells.Item(132, 8).Value = 10898243
$ws.Cells.Item(132, 9).Value = 13528645
$ws.Cells.Item(132, 10).Value = 863.1429000000001
$ws.Cells.Item(132, 11).Value = 40585935
$ws.Cells.Item(132, 12).Value = 2589.4287
$ws.Cells.Item(132, 13).Value = -40583405
$ws.Cells.Item(132, 14).Value = -7649.4287
